$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 76923190  # H33
$ws.Cells.Item(33, 9).Value = 125.083336  # I33
$ws.Cells.Item(33, 11).Value = 125.083336  # K33
$ws.Cells.Item(33, 13).Value = 103.916664  # M33
$ws.Cells.Item(40, 8).Value = 2324.2  # H40
$ws.Cells.Item(40, 9).Value = 2249.111  # I40
$ws.Cells.Item(40, 11).Value = 2249.111  # K40
$ws.Cells.Item(40, 13).Value = -2074.111  # M40
$ws.Cells.Item(43, 8).Value = 10329.667  # H43
$ws.Cells.Item(43, 9).Value = 16499.5  # I43
$ws.Cells.Item(43, 10).Value = 7244.75  # J43
$ws.Cells.Item(43, 11).Value = 16499.5  # K43
$ws.Cells.Item(43, 12).Value = 7244.75  # L43
$ws.Cells.Item(43, 13).Value = -16430.5  # M43
$ws.Cells.Item(43, 14).Value = -7382.75  # N43
$ws.Cells.Item(52, 8).Value = 1623.9166  # H52
$ws.Cells.Item(52, 10).Value = 6999.5  # J52
$ws.Cells.Item(52, 12).Value = 20998.5  # L52
$ws.Cells.Item(52, 14).Value = -21318.5  # N52
$ws.Cells.Item(106, 8).Value = 36416.375  # H106
$ws.Cells.Item(106, 9).Value = 41190.145  # I106
$ws.Cells.Item(106, 11).Value = 41190.145  # K106
$ws.Cells.Item(106, 13).Value = -40559.145  # M106
$ws.Cells.Item(107, 8).Value = 1175.3684  # H107
$ws.Cells.Item(107, 9).Value = 1175.3684  # I107
$ws.Cells.Item(107, 10).Value = 0  # J107
$ws.Cells.Item(107, 11).Value = 1175.3684  # K107
$ws.Cells.Item(107, 12).Value = 0  # L107
$ws.Cells.Item(107, 13).Value = 744.6315999999999  # M107
$ws.Cells.Item(107, 14).ClearContents()  # N107
$ws.Cells.Item(116, 10).Value = 0  # J116
$ws.Cells.Item(116, 12).Value = 0  # L116
$ws.Cells.Item(116, 14).ClearContents()  # N116
$ws.Cells.Item(132, 8).Value = 962.65515  # H132
$ws.Cells.Item(132, 9).Value = 907.7857  # I132
$ws.Cells.Item(132, 11).Value = 2723.3571  # K132
$ws.Cells.Item(132, 13).Value = -193.3571000000002  # M132
$ws.Cells.Item(137, 8).Value = 2311.7727  # H137
$ws.Cells.Item(137, 9).Value = 2870.6924  # I137
$ws.Cells.Item(137, 11).Value = 8612.0772  # K137
$ws.Cells.Item(137, 13).Value = -6062.0772  # M137

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(19, 8).Value = 36124.5  # H19
$ws.Cells.Item(19, 9).Value = 36124.5  # I19
$ws.Cells.Item(19, 11).Value = 36124.5  # K19
$ws.Cells.Item(19, 13).Value = -35895.5  # M19
$ws.Cells.Item(32, 8).Value = 2850.7778  # H32
$ws.Cells.Item(32, 9).Value = 1941.381  # I32
$ws.Cells.Item(32, 11).Value = 1941.381  # K32
$ws.Cells.Item(32, 13).Value = -1654.381  # M32
$ws.Cells.Item(74, 8).Value = 2463.4333  # H74
$ws.Cells.Item(74, 9).Value = 2356.12  # I74
$ws.Cells.Item(74, 10).Value = 3000  # J74
$ws.Cells.Item(74, 11).Value = 2356.12  # K74
$ws.Cells.Item(74, 12).Value = 3000  # L74
$ws.Cells.Item(74, 13).Value = -1482.12  # M74
$ws.Cells.Item(74, 14).Value = -4748  # N74
$ws.Cells.Item(77, 8).Value = 2463.4333  # H77
$ws.Cells.Item(77, 9).Value = 2356.12  # I77
$ws.Cells.Item(77, 10).Value = 3000  # J77
$ws.Cells.Item(77, 11).Value = 11780.6  # K77
$ws.Cells.Item(77, 12).Value = 15000  # L77
$ws.Cells.Item(77, 13).Value = -7412.599999999999  # M77
$ws.Cells.Item(77, 14).Value = -23736  # N77
$ws.Cells.Item(93, 8).Value = 15000  # H93
$ws.Cells.Item(93, 10).Value = 15000  # J93
$ws.Cells.Item(93, 12).Value = 15000  # L93
$ws.Cells.Item(93, 14).Value = -19992  # N93
$ws.Cells.Item(97, 8).Value = 449.68182  # H97
$ws.Cells.Item(97, 9).Value = 243.625  # I97
$ws.Cells.Item(97, 10).Value = 999.1667  # J97
$ws.Cells.Item(97, 11).Value = 243.625  # K97
$ws.Cells.Item(97, 12).Value = 999.1667  # L97
$ws.Cells.Item(97, 13).Value = 252.375  # M97
$ws.Cells.Item(97, 14).Value = -1991.1667  # N97
$ws.Cells.Item(122, 8).Value = 3188.465  # H122
$ws.Cells.Item(122, 9).Value = 2891.861  # I122
$ws.Cells.Item(122, 10).Value = 4713.857  # J122
$ws.Cells.Item(122, 11).Value = 8675.582999999999  # K122
$ws.Cells.Item(122, 12).Value = 14141.571  # L122
$ws.Cells.Item(122, 13).Value = -6225.582999999999  # M122
$ws.Cells.Item(122, 14).Value = -19041.571  # N122
$ws.Cells.Item(123, 8).Value = 61300  # H123
$ws.Cells.Item(123, 10).Value = 61300  # J123
$ws.Cells.Item(123, 12).Value = 61300  # L123
$ws.Cells.Item(123, 14).Value = -71100  # N123

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1177.4445  # H20
$ws.Cells.Item(20, 9).Value = 1233.8948  # I20
$ws.Cells.Item(20, 11).Value = 1233.8948  # K20
$ws.Cells.Item(20, 13).Value = -986.8948  # M20
$ws.Cells.Item(36, 8).Value = 3437.8  # H36
$ws.Cells.Item(36, 9).Value = 3437.8  # I36
$ws.Cells.Item(36, 11).Value = 3437.8  # K36
$ws.Cells.Item(36, 13).Value = -2903.8  # M36
$ws.Cells.Item(94, 8).Value = 4808.625  # H94
$ws.Cells.Item(94, 10).Value = 5282  # J94
$ws.Cells.Item(94, 12).Value = 5282  # L94
$ws.Cells.Item(94, 14).Value = -6184  # N94

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1345.7693  # H22
$ws.Cells.Item(22, 9).Value = 362.25  # I22
$ws.Cells.Item(22, 10).Value = 1782.8889  # J22
$ws.Cells.Item(22, 11).Value = 362.25  # K22
$ws.Cells.Item(22, 12).Value = 1782.8889  # L22
$ws.Cells.Item(22, 13).Value = -12.25  # M22
$ws.Cells.Item(22, 14).Value = -2482.8889  # N22
$ws.Cells.Item(31, 8).Value = 1294.4286  # H31
$ws.Cells.Item(31, 9).Value = 1332.8462  # I31
$ws.Cells.Item(31, 10).Value = 795  # J31
$ws.Cells.Item(31, 11).Value = 1332.8462  # K31
$ws.Cells.Item(31, 12).Value = 795  # L31
$ws.Cells.Item(31, 13).Value = -1037.8462  # M31
$ws.Cells.Item(31, 14).Value = -1385  # N31
$ws.Cells.Item(34, 8).Value = 1294.4286  # H34
$ws.Cells.Item(34, 9).Value = 1332.8462  # I34
$ws.Cells.Item(34, 10).Value = 795  # J34
$ws.Cells.Item(34, 11).Value = 1332.8462  # K34
$ws.Cells.Item(34, 12).Value = 795  # L34
$ws.Cells.Item(34, 13).Value = -1130.8462  # M34
$ws.Cells.Item(34, 14).Value = -1199  # N34

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1100.8  # H5
$ws.Cells.Item(5, 9).Value = 1156.5555  # I5
$ws.Cells.Item(5, 11).Value = 3469.6665  # K5
$ws.Cells.Item(5, 13).Value = -3357.6665  # M5
$ws.Cells.Item(11, 8).Value = 314.83334  # H11
$ws.Cells.Item(11, 9).Value = 397.25  # I11
$ws.Cells.Item(11, 11).Value = 1191.75  # K11
$ws.Cells.Item(11, 13).Value = -1051.75  # M11
$ws.Cells.Item(38, 8).Value = 139.75  # H38
$ws.Cells.Item(38, 10).Value = 194  # J38
$ws.Cells.Item(38, 12).Value = 582  # L38
$ws.Cells.Item(38, 14).Value = -1276  # N38
$ws.Cells.Item(39, 8).Value = 3494.182  # H39
$ws.Cells.Item(39, 9).Value = 400  # I39
$ws.Cells.Item(39, 10).Value = 3803.6  # J39
$ws.Cells.Item(39, 11).Value = 1200  # K39
$ws.Cells.Item(39, 12).Value = 11410.8  # L39
$ws.Cells.Item(39, 13).Value = -906  # M39
$ws.Cells.Item(39, 14).Value = -11998.8  # N39
$ws.Cells.Item(121, 8).Value = 1833.9524  # H121
$ws.Cells.Item(121, 10).Value = 2062.611  # J121
$ws.Cells.Item(121, 12).Value = 6187.833  # L121
$ws.Cells.Item(121, 14).Value = -8807.832999999999  # N121
$ws.Cells.Item(131, 8).Value = 1506.6666  # H131
$ws.Cells.Item(131, 10).Value = 1899.091  # J131
$ws.Cells.Item(131, 12).Value = 5697.272999999999  # L131
$ws.Cells.Item(131, 14).Value = -15777.273  # N131
$ws.Cells.Item(135, 8).Value = 1100.8  # H135
$ws.Cells.Item(135, 9).Value = 1156.5555  # I135
$ws.Cells.Item(135, 11).Value = 10408.9995  # K135
$ws.Cells.Item(135, 13).Value = -7873.9995  # M135

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 225  # H2
$ws.Cells.Item(2, 9).Value = 187.84616  # I2
$ws.Cells.Item(2, 11).Value = 187.84616  # K2
$ws.Cells.Item(2, 13).Value = -74.84616  # M2
$ws.Cells.Item(21, 8).Value = 37749.5  # H21
$ws.Cells.Item(30, 8).Value = 37749.5  # H30
$ws.Cells.Item(70, 8).Value = 28222.629  # H70
$ws.Cells.Item(70, 9).Value = 37612.06  # I70
$ws.Cells.Item(70, 11).Value = 37612.06  # K70
$ws.Cells.Item(70, 13).Value = -37342.06  # M70
$ws.Cells.Item(73, 8).Value = 28222.629  # H73
$ws.Cells.Item(73, 9).Value = 37612.06  # I73
$ws.Cells.Item(73, 11).Value = 37612.06  # K73
$ws.Cells.Item(73, 13).Value = -36676.06  # M73
$ws.Cells.Item(102, 8).Value = 2258.4736  # H102
$ws.Cells.Item(102, 9).Value = 2008.6428  # I102
$ws.Cells.Item(102, 10).Value = 2958  # J102
$ws.Cells.Item(102, 11).Value = 2008.6428  # K102
$ws.Cells.Item(102, 12).Value = 2958  # L102
$ws.Cells.Item(102, 13).Value = -386.6428000000001  # M102
$ws.Cells.Item(102, 14).Value = -6202  # N102
$ws.Cells.Item(132, 8).Value = 2885.5454  # H132
$ws.Cells.Item(132, 9).Value = 2416.7778  # I132
$ws.Cells.Item(132, 11).Value = 7250.3334  # K132
$ws.Cells.Item(132, 13).Value = -4720.3334  # M132
$ws.Cells.Item(134, 8).Value = 0  # H134
$ws.Cells.Item(134, 10).Value = 0  # J134
$ws.Cells.Item(134, 12).Value = 0  # L134
$ws.Cells.Item(134, 14).ClearContents()  # N134

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1208.875  # H16
$ws.Cells.Item(16, 9).Value = 1014.3333  # I16
$ws.Cells.Item(16, 10).Value = 1533.1111  # J16
$ws.Cells.Item(16, 11).Value = 1014.3333  # K16
$ws.Cells.Item(16, 12).Value = 1533.1111  # L16
$ws.Cells.Item(16, 13).Value = -844.3333  # M16
$ws.Cells.Item(16, 14).Value = -1873.1111  # N16
$ws.Cells.Item(23, 8).Value = 0  # H23
$ws.Cells.Item(23, 9).Value = 0  # I23
$ws.Cells.Item(23, 11).Value = 0  # K23
$ws.Cells.Item(23, 13).ClearContents()  # M23
$ws.Cells.Item(122, 8).Value = 6678.9644  # H122
$ws.Cells.Item(122, 9).Value = 5000.9  # I122
$ws.Cells.Item(122, 11).Value = 15002.7  # K122
$ws.Cells.Item(122, 13).Value = -12552.7  # M122

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1721.1111  # H132
$ws.Cells.Item(132, 9).Value = 1721.1111  # I132
$ws.Cells.Item(132, 10).Value = 0  # J132
$ws.Cells.Item(132, 11).Value = 5163.3333  # K132
$ws.Cells.Item(132, 12).Value = 0  # L132
$ws.Cells.Item(132, 13).Value = -2633.3333  # M132
$ws.Cells.Item(132, 14).ClearContents()  # N132
